$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 322. This shifts the existing rows 322..382
# down to 323..383, preserving all of their data/formatting.
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new record's data.
$ws.Cells.Item(322, 1).Value2 = 7
$ws.Cells.Item(322, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(322, 3).Value2 = "Ñuble"
$ws.Cells.Item(322, 4).Value2 = 44694
$ws.Cells.Item(322, 5).Value2 = 16
$ws.Cells.Item(322, 6).Value2 = 100114014
$ws.Cells.Item(322, 7).Value2 = "Betarraga"
$ws.Cells.Item(322, 8).Value2 = "Sin especificar"
$ws.Cells.Item(322, 9).Value2 = "Segunda"
$ws.Cells.Item(322, 10).Value2 = 150
$ws.Cells.Item(322, 11).Value2 = 600
$ws.Cells.Item(322, 12).Value2 = 600
$ws.Cells.Item(322, 13).Value2 = 600
$ws.Cells.Item(322, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(322, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(322, 16).Value2 = 120
$ws.Cells.Item(322, 17).Value2 = 5
$ws.Cells.Item(322, 18).Value2 = "Hortaliza"
